$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Ruiz"
$ws.Range("B6").Value = "Alzate"
$ws.Range("C6").Value = 32434233423
$ws.Range("D6").Value = 123124314
$ws.Range("E6").Value = "Ruiz@gmail.com"

$ws.Range("A7").Value = "Julia"
$ws.Range("B7").Value = "ana"
$ws.Range("C7").Value = 424324234
$ws.Range("D7").Value = 2342342
$ws.Range("E7").Value = "Juli@hotmail.com"
